$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 61-68: content is cyclically shifted up by one row.
# (old row62 -> row61, old row63 -> row62, ..., old row68 -> row67, old row61 -> row68)
# Every cell A:AY for rows 61-68 is set explicitly to its target value below.
# NumberFormat is forced to text ("@") before assigning string values so that
# date-like ("2019-10-21") or numeric-like ("2") text keeps its original text type
# instead of being auto-coerced into a date serial / number by Excel.

$ws.Range("A61").Value = 81745276
$ws.Range("B61").Value = 5135
$ws.Range("C61").NumberFormat = "@"
$ws.Range("C61").Value = "Ovaliderad"
$ws.Range("D61").NumberFormat = "@"
$ws.Range("D61").Value = "LC"
$ws.Range("E61").Value = 105930
$ws.Range("F61").NumberFormat = "@"
$ws.Range("F61").Value = "Vågbandad barkbock"
$ws.Range("G61").NumberFormat = "@"
$ws.Range("G61").Value = "Semanotus undatus"
$ws.Range("H61").NumberFormat = "@"
$ws.Range("H61").Value = "(Linnaeus, 1758)"
$ws.Range("I61").NumberFormat = "@"
$ws.Range("I61").Value = ""
$ws.Range("J61").NumberFormat = "@"
$ws.Range("J61").Value = ""
$ws.Range("K61").NumberFormat = "@"
$ws.Range("K61").Value = ""
$ws.Range("L61").NumberFormat = "@"
$ws.Range("L61").Value = ""
$ws.Range("M61").NumberFormat = "@"
$ws.Range("M61").Value = "äldre gnagspår"
$ws.Range("N61").NumberFormat = "@"
$ws.Range("N61").Value = ""
$ws.Range("O61").NumberFormat = "@"
$ws.Range("O61").Value = ""
$ws.Range("P61").NumberFormat = "@"
$ws.Range("P61").Value = "Valåsen, Vrm"
$ws.Range("Q61").Value = 318237.8515678492
$ws.Range("R61").Value = 6591495.875899103
$ws.Range("S61").Value = 10
$ws.Range("T61").NumberFormat = "@"
$ws.Range("T61").Value = "Värmland"
$ws.Range("U61").NumberFormat = "@"
$ws.Range("U61").Value = "Årjäng"
$ws.Range("V61").NumberFormat = "@"
$ws.Range("V61").Value = "Värmland"
$ws.Range("W61").NumberFormat = "@"
$ws.Range("W61").Value = "Västra Fågelvik"
$ws.Range("X61").NumberFormat = "@"
$ws.Range("X61").Value = ""
$ws.Range("Y61").NumberFormat = "@"
$ws.Range("Y61").Value = "2019-10-21"
$ws.Range("Z61").NumberFormat = "@"
$ws.Range("Z61").Value = "00:00"
$ws.Range("AA61").NumberFormat = "@"
$ws.Range("AA61").Value = "2019-10-21"
$ws.Range("AB61").NumberFormat = "@"
$ws.Range("AB61").Value = "00:00"
$ws.Range("AC61").NumberFormat = "@"
$ws.Range("AC61").Value = ""
$ws.Range("AD61").Value = $false
$ws.Range("AE61").Value = $false
$ws.Range("AF61").NumberFormat = "@"
$ws.Range("AF61").Value = ""
$ws.Range("AG61").Value = $false
$ws.Range("AH61").NumberFormat = "@"
$ws.Range("AH61").Value = ""
$ws.Range("AI61").NumberFormat = "@"
$ws.Range("AI61").Value = ""
$ws.Range("AJ61").NumberFormat = "@"
$ws.Range("AJ61").Value = ""
$ws.Range("AK61").NumberFormat = "@"
$ws.Range("AK61").Value = ""
$ws.Range("AL61").NumberFormat = "@"
$ws.Range("AL61").Value = ""
$ws.Range("AM61").NumberFormat = "@"
$ws.Range("AM61").Value = ""
$ws.Range("AN61").NumberFormat = "@"
$ws.Range("AN61").Value = ""
$ws.Range("AO61").NumberFormat = "@"
$ws.Range("AO61").Value = ""
$ws.Range("AP61").NumberFormat = "@"
$ws.Range("AP61").Value = ""
$ws.Range("AQ61").NumberFormat = "@"
$ws.Range("AQ61").Value = ""
$ws.Range("AR61").NumberFormat = "@"
$ws.Range("AR61").Value = ""
$ws.Range("AS61").NumberFormat = "@"
$ws.Range("AS61").Value = ""
$ws.Range("AT61").NumberFormat = "@"
$ws.Range("AT61").Value = ""
$ws.Range("AU61").NumberFormat = "@"
$ws.Range("AU61").Value = ""
$ws.Range("AV61").NumberFormat = "@"
$ws.Range("AV61").Value = ""
$ws.Range("AW61").NumberFormat = "@"
$ws.Range("AW61").Value = "Jan Rees"
$ws.Range("AX61").NumberFormat = "@"
$ws.Range("AX61").Value = "Jan Rees, Roger Gran"
$ws.Range("AY61").NumberFormat = "@"
$ws.Range("AY61").Value = ""
$ws.Range("A62").Value = 78458025
$ws.Range("B62").Value = 5135
$ws.Range("C62").NumberFormat = "@"
$ws.Range("C62").Value = "Ovaliderad"
$ws.Range("D62").NumberFormat = "@"
$ws.Range("D62").Value = "LC"
$ws.Range("E62").Value = 105930
$ws.Range("F62").NumberFormat = "@"
$ws.Range("F62").Value = "Vågbandad barkbock"
$ws.Range("G62").NumberFormat = "@"
$ws.Range("G62").Value = "Semanotus undatus"
$ws.Range("H62").NumberFormat = "@"
$ws.Range("H62").Value = "(Linnaeus, 1758)"
$ws.Range("I62").NumberFormat = "@"
$ws.Range("I62").Value = ""
$ws.Range("J62").NumberFormat = "@"
$ws.Range("J62").Value = ""
$ws.Range("K62").NumberFormat = "@"
$ws.Range("K62").Value = ""
$ws.Range("L62").NumberFormat = "@"
$ws.Range("L62").Value = ""
$ws.Range("M62").NumberFormat = "@"
$ws.Range("M62").Value = "färska gnagspår"
$ws.Range("N62").NumberFormat = "@"
$ws.Range("N62").Value = ""
$ws.Range("O62").NumberFormat = "@"
$ws.Range("O62").Value = ""
$ws.Range("P62").NumberFormat = "@"
$ws.Range("P62").Value = "Sydväst Valåstjärn, Vrm"
$ws.Range("Q62").Value = 318086.1680672934
$ws.Range("R62").Value = 6591484.330946958
$ws.Range("S62").Value = 10
$ws.Range("T62").NumberFormat = "@"
$ws.Range("T62").Value = "Värmland"
$ws.Range("U62").NumberFormat = "@"
$ws.Range("U62").Value = "Årjäng"
$ws.Range("V62").NumberFormat = "@"
$ws.Range("V62").Value = "Värmland"
$ws.Range("W62").NumberFormat = "@"
$ws.Range("W62").Value = "Västra Fågelvik"
$ws.Range("X62").NumberFormat = "@"
$ws.Range("X62").Value = ""
$ws.Range("Y62").NumberFormat = "@"
$ws.Range("Y62").Value = "2019-06-05"
$ws.Range("Z62").NumberFormat = "@"
$ws.Range("Z62").Value = "00:00"
$ws.Range("AA62").NumberFormat = "@"
$ws.Range("AA62").Value = "2019-06-05"
$ws.Range("AB62").NumberFormat = "@"
$ws.Range("AB62").Value = "00:00"
$ws.Range("AC62").NumberFormat = "@"
$ws.Range("AC62").Value = ""
$ws.Range("AD62").Value = $false
$ws.Range("AE62").Value = $false
$ws.Range("AF62").NumberFormat = "@"
$ws.Range("AF62").Value = ""
$ws.Range("AG62").Value = $false
$ws.Range("AH62").NumberFormat = "@"
$ws.Range("AH62").Value = ""
$ws.Range("AI62").NumberFormat = "@"
$ws.Range("AI62").Value = ""
$ws.Range("AJ62").NumberFormat = "@"
$ws.Range("AJ62").Value = ""
$ws.Range("AK62").NumberFormat = "@"
$ws.Range("AK62").Value = ""
$ws.Range("AL62").NumberFormat = "@"
$ws.Range("AL62").Value = ""
$ws.Range("AM62").NumberFormat = "@"
$ws.Range("AM62").Value = ""
$ws.Range("AN62").NumberFormat = "@"
$ws.Range("AN62").Value = ""
$ws.Range("AO62").NumberFormat = "@"
$ws.Range("AO62").Value = ""
$ws.Range("AP62").NumberFormat = "@"
$ws.Range("AP62").Value = ""
$ws.Range("AQ62").NumberFormat = "@"
$ws.Range("AQ62").Value = ""
$ws.Range("AR62").NumberFormat = "@"
$ws.Range("AR62").Value = ""
$ws.Range("AS62").NumberFormat = "@"
$ws.Range("AS62").Value = ""
$ws.Range("AT62").NumberFormat = "@"
$ws.Range("AT62").Value = ""
$ws.Range("AU62").NumberFormat = "@"
$ws.Range("AU62").Value = ""
$ws.Range("AV62").NumberFormat = "@"
$ws.Range("AV62").Value = ""
$ws.Range("AW62").NumberFormat = "@"
$ws.Range("AW62").Value = "Roger Gran"
$ws.Range("AX62").NumberFormat = "@"
$ws.Range("AX62").Value = "Roger Gran"
$ws.Range("AY62").NumberFormat = "@"
$ws.Range("AY62").Value = ""
$ws.Range("A63").Value = 106628007
$ws.Range("B63").Value = 89997
$ws.Range("C63").NumberFormat = "@"
$ws.Range("C63").Value = "Ovaliderad"
$ws.Range("D63").NumberFormat = "@"
$ws.Range("D63").Value = "NT"
$ws.Range("E63").Value = 5454
$ws.Range("F63").NumberFormat = "@"
$ws.Range("F63").Value = "Hornvaxskinn"
$ws.Range("G63").NumberFormat = "@"
$ws.Range("G63").Value = "Crustoderma corneum"
$ws.Range("H63").NumberFormat = "@"
$ws.Range("H63").Value = "(Bourdot & Galzin) Nakasone"
$ws.Range("I63").NumberFormat = "@"
$ws.Range("I63").Value = ""
$ws.Range("J63").NumberFormat = "@"
$ws.Range("J63").Value = "fruktkroppar"
$ws.Range("K63").NumberFormat = "@"
$ws.Range("K63").Value = ""
$ws.Range("L63").NumberFormat = "@"
$ws.Range("L63").Value = ""
$ws.Range("M63").NumberFormat = "@"
$ws.Range("M63").Value = ""
$ws.Range("N63").NumberFormat = "@"
$ws.Range("N63").Value = ""
$ws.Range("O63").NumberFormat = "@"
$ws.Range("O63").Value = ""
$ws.Range("P63").NumberFormat = "@"
$ws.Range("P63").Value = "Söder om Valåsen, Vrm"
$ws.Range("Q63").Value = 317838.0016589427
$ws.Range("R63").Value = 6591494.270068384
$ws.Range("S63").Value = 10
$ws.Range("T63").NumberFormat = "@"
$ws.Range("T63").Value = "Värmland"
$ws.Range("U63").NumberFormat = "@"
$ws.Range("U63").Value = "Årjäng"
$ws.Range("V63").NumberFormat = "@"
$ws.Range("V63").Value = "Värmland"
$ws.Range("W63").NumberFormat = "@"
$ws.Range("W63").Value = "Västra Fågelvik"
$ws.Range("X63").NumberFormat = "@"
$ws.Range("X63").Value = ""
$ws.Range("Y63").NumberFormat = "@"
$ws.Range("Y63").Value = "2023-02-11"
$ws.Range("Z63").NumberFormat = "@"
$ws.Range("Z63").Value = "00:00"
$ws.Range("AA63").NumberFormat = "@"
$ws.Range("AA63").Value = "2023-02-11"
$ws.Range("AB63").NumberFormat = "@"
$ws.Range("AB63").Value = "00:00"
$ws.Range("AC63").NumberFormat = "@"
$ws.Range("AC63").Value = ""
$ws.Range("AD63").Value = $false
$ws.Range("AE63").Value = $false
$ws.Range("AF63").NumberFormat = "@"
$ws.Range("AF63").Value = ""
$ws.Range("AG63").Value = $false
$ws.Range("AH63").NumberFormat = "@"
$ws.Range("AH63").Value = ""
$ws.Range("AI63").NumberFormat = "@"
$ws.Range("AI63").Value = "Barrblandskog"
$ws.Range("AJ63").NumberFormat = "@"
$ws.Range("AJ63").Value = "tall"
$ws.Range("AK63").NumberFormat = "@"
$ws.Range("AK63").Value = "Pinus sylvestris"
$ws.Range("AL63").NumberFormat = "@"
$ws.Range("AL63").Value = ""
$ws.Range("AM63").NumberFormat = "@"
$ws.Range("AM63").Value = ""
$ws.Range("AN63").NumberFormat = "@"
$ws.Range("AN63").Value = ""
$ws.Range("AO63").NumberFormat = "@"
$ws.Range("AO63").Value = "Låga # Pinus sylvestris"
$ws.Range("AP63").NumberFormat = "@"
$ws.Range("AP63").Value = ""
$ws.Range("AQ63").NumberFormat = "@"
$ws.Range("AQ63").Value = ""
$ws.Range("AR63").NumberFormat = "@"
$ws.Range("AR63").Value = ""
$ws.Range("AS63").NumberFormat = "@"
$ws.Range("AS63").Value = ""
$ws.Range("AT63").NumberFormat = "@"
$ws.Range("AT63").Value = ""
$ws.Range("AU63").NumberFormat = "@"
$ws.Range("AU63").Value = ""
$ws.Range("AV63").NumberFormat = "@"
$ws.Range("AV63").Value = ""
$ws.Range("AW63").NumberFormat = "@"
$ws.Range("AW63").Value = "Bård E. Andersen"
$ws.Range("AX63").NumberFormat = "@"
$ws.Range("AX63").Value = "Bård E. Andersen"
$ws.Range("AY63").NumberFormat = "@"
$ws.Range("AY63").Value = ""
$ws.Range("A64").Value = 106627208
$ws.Range("B64").Value = 89997
$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").Value = "Ovaliderad"
$ws.Range("D64").NumberFormat = "@"
$ws.Range("D64").Value = "NT"
$ws.Range("E64").Value = 5454
$ws.Range("F64").NumberFormat = "@"
$ws.Range("F64").Value = "Hornvaxskinn"
$ws.Range("G64").NumberFormat = "@"
$ws.Range("G64").Value = "Crustoderma corneum"
$ws.Range("H64").NumberFormat = "@"
$ws.Range("H64").Value = "(Bourdot & Galzin) Nakasone"
$ws.Range("I64").NumberFormat = "@"
$ws.Range("I64").Value = ""
$ws.Range("J64").NumberFormat = "@"
$ws.Range("J64").Value = "fruktkroppar"
$ws.Range("K64").NumberFormat = "@"
$ws.Range("K64").Value = ""
$ws.Range("L64").NumberFormat = "@"
$ws.Range("L64").Value = ""
$ws.Range("M64").NumberFormat = "@"
$ws.Range("M64").Value = ""
$ws.Range("N64").NumberFormat = "@"
$ws.Range("N64").Value = ""
$ws.Range("O64").NumberFormat = "@"
$ws.Range("O64").Value = ""
$ws.Range("P64").NumberFormat = "@"
$ws.Range("P64").Value = "Vängsfjället, Vrm"
$ws.Range("Q64").Value = 317636.1213810217
$ws.Range("R64").Value = 6591247.138022185
$ws.Range("S64").Value = 10
$ws.Range("T64").NumberFormat = "@"
$ws.Range("T64").Value = "Värmland"
$ws.Range("U64").NumberFormat = "@"
$ws.Range("U64").Value = "Årjäng"
$ws.Range("V64").NumberFormat = "@"
$ws.Range("V64").Value = "Värmland"
$ws.Range("W64").NumberFormat = "@"
$ws.Range("W64").Value = "Västra Fågelvik"
$ws.Range("X64").NumberFormat = "@"
$ws.Range("X64").Value = ""
$ws.Range("Y64").NumberFormat = "@"
$ws.Range("Y64").Value = "2023-02-11"
$ws.Range("Z64").NumberFormat = "@"
$ws.Range("Z64").Value = "00:00"
$ws.Range("AA64").NumberFormat = "@"
$ws.Range("AA64").Value = "2023-02-11"
$ws.Range("AB64").NumberFormat = "@"
$ws.Range("AB64").Value = "00:00"
$ws.Range("AC64").NumberFormat = "@"
$ws.Range("AC64").Value = ""
$ws.Range("AD64").Value = $false
$ws.Range("AE64").Value = $false
$ws.Range("AF64").NumberFormat = "@"
$ws.Range("AF64").Value = ""
$ws.Range("AG64").Value = $false
$ws.Range("AH64").NumberFormat = "@"
$ws.Range("AH64").Value = ""
$ws.Range("AI64").NumberFormat = "@"
$ws.Range("AI64").Value = ""
$ws.Range("AJ64").NumberFormat = "@"
$ws.Range("AJ64").Value = ""
$ws.Range("AK64").NumberFormat = "@"
$ws.Range("AK64").Value = ""
$ws.Range("AL64").NumberFormat = "@"
$ws.Range("AL64").Value = ""
$ws.Range("AM64").NumberFormat = "@"
$ws.Range("AM64").Value = ""
$ws.Range("AN64").NumberFormat = "@"
$ws.Range("AN64").Value = ""
$ws.Range("AO64").NumberFormat = "@"
$ws.Range("AO64").Value = ""
$ws.Range("AP64").NumberFormat = "@"
$ws.Range("AP64").Value = ""
$ws.Range("AQ64").NumberFormat = "@"
$ws.Range("AQ64").Value = ""
$ws.Range("AR64").NumberFormat = "@"
$ws.Range("AR64").Value = ""
$ws.Range("AS64").NumberFormat = "@"
$ws.Range("AS64").Value = ""
$ws.Range("AT64").NumberFormat = "@"
$ws.Range("AT64").Value = ""
$ws.Range("AU64").NumberFormat = "@"
$ws.Range("AU64").Value = ""
$ws.Range("AV64").NumberFormat = "@"
$ws.Range("AV64").Value = ""
$ws.Range("AW64").NumberFormat = "@"
$ws.Range("AW64").Value = "Bård E. Andersen"
$ws.Range("AX64").NumberFormat = "@"
$ws.Range("AX64").Value = "Bård E. Andersen"
$ws.Range("AY64").NumberFormat = "@"
$ws.Range("AY64").Value = ""
$ws.Range("A65").Value = 106627422
$ws.Range("B65").Value = 89997
$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").Value = "Ovaliderad"
$ws.Range("D65").NumberFormat = "@"
$ws.Range("D65").Value = "NT"
$ws.Range("E65").Value = 5454
$ws.Range("F65").NumberFormat = "@"
$ws.Range("F65").Value = "Hornvaxskinn"
$ws.Range("G65").NumberFormat = "@"
$ws.Range("G65").Value = "Crustoderma corneum"
$ws.Range("H65").NumberFormat = "@"
$ws.Range("H65").Value = "(Bourdot & Galzin) Nakasone"
$ws.Range("I65").NumberFormat = "@"
$ws.Range("I65").Value = ""
$ws.Range("J65").NumberFormat = "@"
$ws.Range("J65").Value = "fruktkroppar"
$ws.Range("K65").NumberFormat = "@"
$ws.Range("K65").Value = ""
$ws.Range("L65").NumberFormat = "@"
$ws.Range("L65").Value = ""
$ws.Range("M65").NumberFormat = "@"
$ws.Range("M65").Value = ""
$ws.Range("N65").NumberFormat = "@"
$ws.Range("N65").Value = ""
$ws.Range("O65").NumberFormat = "@"
$ws.Range("O65").Value = ""
$ws.Range("P65").NumberFormat = "@"
$ws.Range("P65").Value = "Söder om Tutjärnen, Vrm"
$ws.Range("Q65").Value = 317574.1696661118
$ws.Range("R65").Value = 6591254.211470196
$ws.Range("S65").Value = 10
$ws.Range("T65").NumberFormat = "@"
$ws.Range("T65").Value = "Värmland"
$ws.Range("U65").NumberFormat = "@"
$ws.Range("U65").Value = "Årjäng"
$ws.Range("V65").NumberFormat = "@"
$ws.Range("V65").Value = "Värmland"
$ws.Range("W65").NumberFormat = "@"
$ws.Range("W65").Value = "Västra Fågelvik"
$ws.Range("X65").NumberFormat = "@"
$ws.Range("X65").Value = ""
$ws.Range("Y65").NumberFormat = "@"
$ws.Range("Y65").Value = "2023-02-11"
$ws.Range("Z65").NumberFormat = "@"
$ws.Range("Z65").Value = "00:00"
$ws.Range("AA65").NumberFormat = "@"
$ws.Range("AA65").Value = "2023-02-11"
$ws.Range("AB65").NumberFormat = "@"
$ws.Range("AB65").Value = "00:00"
$ws.Range("AC65").NumberFormat = "@"
$ws.Range("AC65").Value = ""
$ws.Range("AD65").Value = $false
$ws.Range("AE65").Value = $false
$ws.Range("AF65").NumberFormat = "@"
$ws.Range("AF65").Value = ""
$ws.Range("AG65").Value = $false
$ws.Range("AH65").NumberFormat = "@"
$ws.Range("AH65").Value = ""
$ws.Range("AI65").NumberFormat = "@"
$ws.Range("AI65").Value = "Barrblandskog"
$ws.Range("AJ65").NumberFormat = "@"
$ws.Range("AJ65").Value = "tall"
$ws.Range("AK65").NumberFormat = "@"
$ws.Range("AK65").Value = "Pinus sylvestris"
$ws.Range("AL65").NumberFormat = "@"
$ws.Range("AL65").Value = ""
$ws.Range("AM65").NumberFormat = "@"
$ws.Range("AM65").Value = ""
$ws.Range("AN65").NumberFormat = "@"
$ws.Range("AN65").Value = ""
$ws.Range("AO65").NumberFormat = "@"
$ws.Range("AO65").Value = "Låga # Pinus sylvestris"
$ws.Range("AP65").NumberFormat = "@"
$ws.Range("AP65").Value = ""
$ws.Range("AQ65").NumberFormat = "@"
$ws.Range("AQ65").Value = ""
$ws.Range("AR65").NumberFormat = "@"
$ws.Range("AR65").Value = ""
$ws.Range("AS65").NumberFormat = "@"
$ws.Range("AS65").Value = ""
$ws.Range("AT65").NumberFormat = "@"
$ws.Range("AT65").Value = ""
$ws.Range("AU65").NumberFormat = "@"
$ws.Range("AU65").Value = ""
$ws.Range("AV65").NumberFormat = "@"
$ws.Range("AV65").Value = ""
$ws.Range("AW65").NumberFormat = "@"
$ws.Range("AW65").Value = "Bård E. Andersen"
$ws.Range("AX65").NumberFormat = "@"
$ws.Range("AX65").Value = "Bård E. Andersen"
$ws.Range("AY65").NumberFormat = "@"
$ws.Range("AY65").Value = ""
$ws.Range("A66").Value = 108130364
$ws.Range("B66").Value = 56540
$ws.Range("C66").NumberFormat = "@"
$ws.Range("C66").Value = "Ovaliderad"
$ws.Range("D66").NumberFormat = "@"
$ws.Range("D66").Value = "NT"
$ws.Range("E66").Value = 103021
$ws.Range("F66").NumberFormat = "@"
$ws.Range("F66").Value = "Talltita"
$ws.Range("G66").NumberFormat = "@"
$ws.Range("G66").Value = "Poecile montanus"
$ws.Range("H66").NumberFormat = "@"
$ws.Range("H66").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I66").NumberFormat = "@"
$ws.Range("I66").Value = "2"
$ws.Range("J66").NumberFormat = "@"
$ws.Range("J66").Value = ""
$ws.Range("K66").NumberFormat = "@"
$ws.Range("K66").Value = ""
$ws.Range("L66").NumberFormat = "@"
$ws.Range("L66").Value = ""
$ws.Range("M66").NumberFormat = "@"
$ws.Range("M66").Value = "lockläte, övriga läten"
$ws.Range("N66").NumberFormat = "@"
$ws.Range("N66").Value = ""
$ws.Range("O66").NumberFormat = "@"
$ws.Range("O66").Value = ""
$ws.Range("P66").NumberFormat = "@"
$ws.Range("P66").Value = "Tutjärnen, Ö om, Vrm"
$ws.Range("Q66").Value = 317495.3811155346
$ws.Range("R66").Value = 6591187.681347501
$ws.Range("S66").Value = 25
$ws.Range("T66").NumberFormat = "@"
$ws.Range("T66").Value = "Värmland"
$ws.Range("U66").NumberFormat = "@"
$ws.Range("U66").Value = "Årjäng"
$ws.Range("V66").NumberFormat = "@"
$ws.Range("V66").Value = "Värmland"
$ws.Range("W66").NumberFormat = "@"
$ws.Range("W66").Value = "Västra Fågelvik"
$ws.Range("X66").NumberFormat = "@"
$ws.Range("X66").Value = ""
$ws.Range("Y66").NumberFormat = "@"
$ws.Range("Y66").Value = "2023-04-14"
$ws.Range("Z66").NumberFormat = "@"
$ws.Range("Z66").Value = "17:31"
$ws.Range("AA66").NumberFormat = "@"
$ws.Range("AA66").Value = "2023-04-14"
$ws.Range("AB66").NumberFormat = "@"
$ws.Range("AB66").Value = "17:31"
$ws.Range("AC66").NumberFormat = "@"
$ws.Range("AC66").Value = ""
$ws.Range("AD66").Value = $false
$ws.Range("AE66").Value = $false
$ws.Range("AF66").NumberFormat = "@"
$ws.Range("AF66").Value = ""
$ws.Range("AG66").Value = $false
$ws.Range("AH66").NumberFormat = "@"
$ws.Range("AH66").Value = ""
$ws.Range("AI66").NumberFormat = "@"
$ws.Range("AI66").Value = ""
$ws.Range("AJ66").NumberFormat = "@"
$ws.Range("AJ66").Value = ""
$ws.Range("AK66").NumberFormat = "@"
$ws.Range("AK66").Value = ""
$ws.Range("AL66").NumberFormat = "@"
$ws.Range("AL66").Value = ""
$ws.Range("AM66").NumberFormat = "@"
$ws.Range("AM66").Value = ""
$ws.Range("AN66").NumberFormat = "@"
$ws.Range("AN66").Value = ""
$ws.Range("AO66").NumberFormat = "@"
$ws.Range("AO66").Value = ""
$ws.Range("AP66").NumberFormat = "@"
$ws.Range("AP66").Value = ""
$ws.Range("AQ66").NumberFormat = "@"
$ws.Range("AQ66").Value = ""
$ws.Range("AR66").NumberFormat = "@"
$ws.Range("AR66").Value = ""
$ws.Range("AS66").NumberFormat = "@"
$ws.Range("AS66").Value = ""
$ws.Range("AT66").NumberFormat = "@"
$ws.Range("AT66").Value = ""
$ws.Range("AU66").NumberFormat = "@"
$ws.Range("AU66").Value = ""
$ws.Range("AV66").NumberFormat = "@"
$ws.Range("AV66").Value = ""
$ws.Range("AW66").NumberFormat = "@"
$ws.Range("AW66").Value = "Uno Skog"
$ws.Range("AX66").NumberFormat = "@"
$ws.Range("AX66").Value = "Uno Skog"
$ws.Range("AY66").NumberFormat = "@"
$ws.Range("AY66").Value = ""
$ws.Range("A67").Value = 109756385
$ws.Range("B67").Value = 89545
$ws.Range("C67").NumberFormat = "@"
$ws.Range("C67").Value = "Ovaliderad"
$ws.Range("D67").NumberFormat = "@"
$ws.Range("D67").Value = "VU"
$ws.Range("E67").Value = 1503
$ws.Range("F67").NumberFormat = "@"
$ws.Range("F67").Value = "Gräddporing"
$ws.Range("G67").NumberFormat = "@"
$ws.Range("G67").Value = "Sidera lenis"
$ws.Range("H67").NumberFormat = "@"
$ws.Range("H67").Value = "(P.Karst.) Miettinen"
$ws.Range("I67").NumberFormat = "@"
$ws.Range("I67").Value = ""
$ws.Range("J67").NumberFormat = "@"
$ws.Range("J67").Value = "fruktkroppar"
$ws.Range("K67").NumberFormat = "@"
$ws.Range("K67").Value = ""
$ws.Range("L67").NumberFormat = "@"
$ws.Range("L67").Value = ""
$ws.Range("M67").NumberFormat = "@"
$ws.Range("M67").Value = ""
$ws.Range("N67").NumberFormat = "@"
$ws.Range("N67").Value = ""
$ws.Range("O67").NumberFormat = "@"
$ws.Range("O67").Value = ""
$ws.Range("P67").NumberFormat = "@"
$ws.Range("P67").Value = "Söder om Valåstjärnet, Vrm"
$ws.Range("Q67").Value = 318116.6165045868
$ws.Range("R67").Value = 6591681.631153035
$ws.Range("S67").Value = 10
$ws.Range("T67").NumberFormat = "@"
$ws.Range("T67").Value = "Värmland"
$ws.Range("U67").NumberFormat = "@"
$ws.Range("U67").Value = "Årjäng"
$ws.Range("V67").NumberFormat = "@"
$ws.Range("V67").Value = "Värmland"
$ws.Range("W67").NumberFormat = "@"
$ws.Range("W67").Value = "Västra Fågelvik"
$ws.Range("X67").NumberFormat = "@"
$ws.Range("X67").Value = ""
$ws.Range("Y67").NumberFormat = "@"
$ws.Range("Y67").Value = "2023-06-03"
$ws.Range("Z67").NumberFormat = "@"
$ws.Range("Z67").Value = "00:00"
$ws.Range("AA67").NumberFormat = "@"
$ws.Range("AA67").Value = "2023-06-03"
$ws.Range("AB67").NumberFormat = "@"
$ws.Range("AB67").Value = "00:00"
$ws.Range("AC67").NumberFormat = "@"
$ws.Range("AC67").Value = ""
$ws.Range("AD67").Value = $false
$ws.Range("AE67").Value = $false
$ws.Range("AF67").NumberFormat = "@"
$ws.Range("AF67").Value = ""
$ws.Range("AG67").Value = $false
$ws.Range("AH67").NumberFormat = "@"
$ws.Range("AH67").Value = ""
$ws.Range("AI67").NumberFormat = "@"
$ws.Range("AI67").Value = "Barrblandskog"
$ws.Range("AJ67").NumberFormat = "@"
$ws.Range("AJ67").Value = "tall"
$ws.Range("AK67").NumberFormat = "@"
$ws.Range("AK67").Value = "Pinus sylvestris"
$ws.Range("AL67").NumberFormat = "@"
$ws.Range("AL67").Value = ""
$ws.Range("AM67").NumberFormat = "@"
$ws.Range("AM67").Value = ""
$ws.Range("AN67").NumberFormat = "@"
$ws.Range("AN67").Value = ""
$ws.Range("AO67").NumberFormat = "@"
$ws.Range("AO67").Value = "I nisjor i kärnved # Pinus sylvestris"
$ws.Range("AP67").NumberFormat = "@"
$ws.Range("AP67").Value = ""
$ws.Range("AQ67").NumberFormat = "@"
$ws.Range("AQ67").Value = ""
$ws.Range("AR67").NumberFormat = "@"
$ws.Range("AR67").Value = ""
$ws.Range("AS67").NumberFormat = "@"
$ws.Range("AS67").Value = ""
$ws.Range("AT67").NumberFormat = "@"
$ws.Range("AT67").Value = ""
$ws.Range("AU67").NumberFormat = "@"
$ws.Range("AU67").Value = ""
$ws.Range("AV67").NumberFormat = "@"
$ws.Range("AV67").Value = ""
$ws.Range("AW67").NumberFormat = "@"
$ws.Range("AW67").Value = "Bård E. Andersen"
$ws.Range("AX67").NumberFormat = "@"
$ws.Range("AX67").Value = "Bård E. Andersen"
$ws.Range("AY67").NumberFormat = "@"
$ws.Range("AY67").Value = ""
$ws.Range("A68").Value = 102262757
$ws.Range("B68").Value = 96334
$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "Ovaliderad"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = "VU"
$ws.Range("E68").Value = 220787
$ws.Range("F68").NumberFormat = "@"
$ws.Range("F68").Value = "Knärot"
$ws.Range("G68").NumberFormat = "@"
$ws.Range("G68").Value = "Goodyera repens"
$ws.Range("H68").NumberFormat = "@"
$ws.Range("H68").Value = "(L.) R. Br."
$ws.Range("I68").NumberFormat = "@"
$ws.Range("I68").Value = ""
$ws.Range("J68").NumberFormat = "@"
$ws.Range("J68").Value = "plantor/tuvor"
$ws.Range("K68").NumberFormat = "@"
$ws.Range("K68").Value = ""
$ws.Range("L68").NumberFormat = "@"
$ws.Range("L68").Value = ""
$ws.Range("M68").NumberFormat = "@"
$ws.Range("M68").Value = ""
$ws.Range("N68").NumberFormat = "@"
$ws.Range("N68").Value = ""
$ws.Range("O68").NumberFormat = "@"
$ws.Range("O68").Value = ""
$ws.Range("P68").NumberFormat = "@"
$ws.Range("P68").Value = "Söder om Tutjärnen, Vrm"
$ws.Range("Q68").Value = 317262.6606318104
$ws.Range("R68").Value = 6591084.767739117
$ws.Range("S68").Value = 10
$ws.Range("T68").NumberFormat = "@"
$ws.Range("T68").Value = "Värmland"
$ws.Range("U68").NumberFormat = "@"
$ws.Range("U68").Value = "Årjäng"
$ws.Range("V68").NumberFormat = "@"
$ws.Range("V68").Value = "Värmland"
$ws.Range("W68").NumberFormat = "@"
$ws.Range("W68").Value = "Västra Fågelvik"
$ws.Range("X68").NumberFormat = "@"
$ws.Range("X68").Value = ""
$ws.Range("Y68").NumberFormat = "@"
$ws.Range("Y68").Value = "2022-07-14"
$ws.Range("Z68").NumberFormat = "@"
$ws.Range("Z68").Value = "00:00"
$ws.Range("AA68").NumberFormat = "@"
$ws.Range("AA68").Value = "2022-07-14"
$ws.Range("AB68").NumberFormat = "@"
$ws.Range("AB68").Value = "00:00"
$ws.Range("AC68").NumberFormat = "@"
$ws.Range("AC68").Value = ""
$ws.Range("AD68").Value = $false
$ws.Range("AE68").Value = $false
$ws.Range("AF68").NumberFormat = "@"
$ws.Range("AF68").Value = ""
$ws.Range("AG68").Value = $false
$ws.Range("AH68").NumberFormat = "@"
$ws.Range("AH68").Value = ""
$ws.Range("AI68").NumberFormat = "@"
$ws.Range("AI68").Value = "Barrblandskog"
$ws.Range("AJ68").NumberFormat = "@"
$ws.Range("AJ68").Value = ""
$ws.Range("AK68").NumberFormat = "@"
$ws.Range("AK68").Value = ""
$ws.Range("AL68").NumberFormat = "@"
$ws.Range("AL68").Value = ""
$ws.Range("AM68").NumberFormat = "@"
$ws.Range("AM68").Value = ""
$ws.Range("AN68").NumberFormat = "@"
$ws.Range("AN68").Value = ""
$ws.Range("AO68").NumberFormat = "@"
$ws.Range("AO68").Value = ""
$ws.Range("AP68").NumberFormat = "@"
$ws.Range("AP68").Value = ""
$ws.Range("AQ68").NumberFormat = "@"
$ws.Range("AQ68").Value = ""
$ws.Range("AR68").NumberFormat = "@"
$ws.Range("AR68").Value = ""
$ws.Range("AS68").NumberFormat = "@"
$ws.Range("AS68").Value = ""
$ws.Range("AT68").NumberFormat = "@"
$ws.Range("AT68").Value = ""
$ws.Range("AU68").NumberFormat = "@"
$ws.Range("AU68").Value = ""
$ws.Range("AV68").NumberFormat = "@"
$ws.Range("AV68").Value = ""
$ws.Range("AW68").NumberFormat = "@"
$ws.Range("AW68").Value = "Bård E. Andersen"
$ws.Range("AX68").NumberFormat = "@"
$ws.Range("AX68").Value = "Bård E. Andersen"
$ws.Range("AY68").NumberFormat = "@"
$ws.Range("AY68").Value = ""
